$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($CellRef, $Text) {
    $rng = $ws.Range($CellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $Text
    $rng.Style = "Normal"
}

Set-TextValue "A2" "2025/11/12"
Set-TextValue "B2" "4.57"

Set-TextValue "A8" "2025/11/12"

Set-TextValue "A14" "2025/11/12"

Set-TextValue "A20" "2025/11/12"
Set-TextValue "B20" "12.01"

Set-TextValue "A26" "2025/11/12"
Set-TextValue "B26" "9.80"

Set-TextValue "A32" "2025/11/12"
Set-TextValue "B32" "24.66"

Set-TextValue "A38" "2025/11/12"

Set-TextValue "A44" "2025/11/12"
Set-TextValue "B44" "11.13"

Set-TextValue "A50" "2025/11/12"

Set-TextValue "A56" "2025/11/12"
Set-TextValue "B56" "34.72"

Set-TextValue "A62" "2025/11/12"
Set-TextValue "B62" "11.47"

Set-TextValue "A68" "2025/11/12"
Set-TextValue "B68" "12.83"

Set-TextValue "A74" "2025/11/12"
Set-TextValue "B74" "15.78"
